$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date
$ws.Name = "Through 2022-05-04"

# Update the header label in I1 (shared string "2022 (through 05-03)" -> "2022 (through 05-04)")
$ws.Range("I1").Value = "2022 (through 05-04)"

# Update May (row 6) total for 2022 column (I)
$ws.Range("I6").Value = 12

# Update the grand Total row (row 14) for 2022 column (I)
$ws.Range("I14").Value = 563
